$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting from the last existing data row (row 16) down to the new rows (17-26)
# so the new rows visually match the existing ones (bold/centered/bordered model-name column).
$fmtSrc = $ws.Range("A16:I16")
$fmtSrc.Copy()
$fmtDst = $ws.Range("A17:I26")
$fmtDst.PasteSpecial(-4122)

# Now populate all data rows (2-26) with the refreshed metrics.
$ws.Cells.Item(2, 1).Value = "model_11_8_0"
$ws.Cells.Item(2, 2).Value = 0.7706758378222208
$ws.Cells.Item(2, 3).Value = 0.8724486424972466
$ws.Cells.Item(2, 4).Value = 0.9010330015449755
$ws.Cells.Item(2, 5).Value = 0.9080584101477194
$ws.Cells.Item(2, 6).Value = 0.253794252872467
$ws.Cells.Item(2, 7).Value = 0.1137137785553932
$ws.Cells.Item(2, 8).Value = 0.06712868809700012
$ws.Cells.Item(2, 9).Value = 0.09179136902093887

$ws.Cells.Item(3, 1).Value = "model_11_8_1"
$ws.Cells.Item(3, 2).Value = 0.813224769946297
$ws.Cells.Item(3, 3).Value = 0.8971339053674581
$ws.Cells.Item(3, 4).Value = 0.784003938405666
$ws.Cells.Item(3, 5).Value = 0.8823118856099572
$ws.Cells.Item(3, 6).Value = 0.2067051380872726
$ws.Cells.Item(3, 7).Value = 0.0917065292596817
$ws.Cells.Item(3, 8).Value = 0.1465087682008743
$ws.Cells.Item(3, 9).Value = 0.1174958348274231

$ws.Cells.Item(4, 1).Value = "model_11_8_3"
$ws.Cells.Item(4, 2).Value = 0.8311366140648514
$ws.Cells.Item(4, 3).Value = 0.4761355749495575
$ws.Cells.Item(4, 4).Value = 0.7431662964929988
$ws.Cells.Item(4, 5).Value = 0.6702282181968899
$ws.Cells.Item(4, 6).Value = 0.1868820041418076
$ws.Cells.Item(4, 7).Value = 0.4670322835445404
$ws.Cells.Item(4, 8).Value = 0.1742086708545685
$ws.Cells.Item(4, 9).Value = 0.3292329609394073

$ws.Cells.Item(5, 1).Value = "model_11_8_7"
$ws.Cells.Item(5, 2).Value = 0.834240579395943
$ws.Cells.Item(5, 3).Value = 0.2169567551925661
$ws.Cells.Item(5, 4).Value = 0.7830169894935455
$ws.Cells.Item(5, 5).Value = 0.5604424890448865
$ws.Cells.Item(5, 6).Value = 0.183446854352951
$ws.Cells.Item(5, 7).Value = 0.6980937123298645
$ws.Cells.Item(5, 8).Value = 0.1471782028675079
$ws.Cells.Item(5, 9).Value = 0.438839316368103

$ws.Cells.Item(6, 1).Value = "model_11_8_6"
$ws.Cells.Item(6, 2).Value = 0.8355518513734161
$ws.Cells.Item(6, 3).Value = 0.237747955477288
$ws.Cells.Item(6, 4).Value = 0.8118424512304178
$ws.Cells.Item(6, 5).Value = 0.5794875252526752
$ws.Cells.Item(6, 6).Value = 0.1819956451654434
$ws.Cells.Item(6, 7).Value = 0.6795580983161926
$ws.Cells.Item(6, 8).Value = 0.1276260763406754
$ws.Cells.Item(6, 9).Value = 0.4198254346847534

$ws.Cells.Item(7, 1).Value = "model_11_8_9"
$ws.Cells.Item(7, 2).Value = 0.8371350189026706
$ws.Cells.Item(7, 3).Value = 0.2203518556329263
$ws.Cells.Item(7, 4).Value = 0.754621198772174
$ws.Cells.Item(7, 5).Value = 0.5529687930252509
$ws.Cells.Item(7, 6).Value = 0.1802435666322708
$ws.Cells.Item(7, 7).Value = 0.6950669288635254
$ws.Cells.Item(7, 8).Value = 0.1664389073848724
$ws.Cells.Item(7, 9).Value = 0.4463008046150208

$ws.Cells.Item(8, 1).Value = "model_11_8_8"
$ws.Cells.Item(8, 2).Value = 0.8397675264115237
$ws.Cells.Item(8, 3).Value = 0.2100078440243506
$ws.Cells.Item(8, 4).Value = 0.8052377241713548
$ws.Cells.Item(8, 5).Value = 0.564261780993
$ws.Cells.Item(8, 6).Value = 0.1773301213979721
$ws.Cells.Item(8, 7).Value = 0.7042887806892395
$ws.Cells.Item(8, 8).Value = 0.1321060210466385
$ws.Cells.Item(8, 9).Value = 0.4350262880325317

$ws.Cells.Item(9, 1).Value = "model_11_8_4"
$ws.Cells.Item(9, 2).Value = 0.842115557859374
$ws.Cells.Item(9, 3).Value = 0.45141739253251
$ws.Cells.Item(9, 4).Value = 0.7608490436593305
$ws.Cells.Item(9, 5).Value = 0.6641962026494825
$ws.Cells.Item(9, 6).Value = 0.1747315675020218
$ws.Cells.Item(9, 7).Value = 0.489068865776062
$ws.Cells.Item(9, 8).Value = 0.1622145920991898
$ws.Cells.Item(9, 9).Value = 0.3352551162242889

$ws.Cells.Item(10, 1).Value = "model_11_8_15"
$ws.Cells.Item(10, 2).Value = 0.8442533732320419
$ws.Cells.Item(10, 3).Value = 0.1600017898959962
$ws.Cells.Item(10, 4).Value = 0.7693375194514502
$ws.Cells.Item(10, 5).Value = 0.5291433579789111
$ws.Cells.Item(10, 6).Value = 0.1723656356334686
$ws.Cells.Item(10, 7).Value = 0.748869776725769
$ws.Cells.Item(10, 8).Value = 0.1564569026231766
$ws.Cells.Item(10, 9).Value = 0.4700873196125031

$ws.Cells.Item(11, 1).Value = "model_11_8_20"
$ws.Cells.Item(11, 2).Value = 0.8450166112351192
$ws.Cells.Item(11, 3).Value = 0.1564648335267668
$ws.Cells.Item(11, 4).Value = 0.7605056043541292
$ws.Cells.Item(11, 5).Value = 0.5246475191749105
$ws.Cells.Item(11, 6).Value = 0.1715209484100342
$ws.Cells.Item(11, 7).Value = 0.7520231604576111
$ws.Cells.Item(11, 8).Value = 0.1624475419521332
$ws.Cells.Item(11, 9).Value = 0.4745758175849915

$ws.Cells.Item(12, 1).Value = "model_11_8_21"
$ws.Cells.Item(12, 2).Value = 0.8455887718395864
$ws.Cells.Item(12, 3).Value = 0.160467067762422
$ws.Cells.Item(12, 4).Value = 0.7599957316748495
$ws.Cells.Item(12, 5).Value = 0.5263766947496542
$ws.Cells.Item(12, 6).Value = 0.1708877384662628
$ws.Cells.Item(12, 7).Value = 0.7484551668167114
$ws.Cells.Item(12, 8).Value = 0.1627933830022812
$ws.Cells.Item(12, 9).Value = 0.4728494584560394

$ws.Cells.Item(13, 1).Value = "model_11_8_23"
$ws.Cells.Item(13, 2).Value = 0.8457439152511388
$ws.Cells.Item(13, 3).Value = 0.1759327225451844
$ws.Cells.Item(13, 4).Value = 0.7372879909976795
$ws.Cells.Item(13, 5).Value = 0.5264278593457408
$ws.Cells.Item(13, 6).Value = 0.1707160323858261
$ws.Cells.Item(13, 7).Value = 0.7346672415733337
$ws.Cells.Item(13, 8).Value = 0.1781958937644958
$ws.Cells.Item(13, 9).Value = 0.4727984070777893

$ws.Cells.Item(14, 1).Value = "model_11_8_17"
$ws.Cells.Item(14, 2).Value = 0.8461596150826338
$ws.Cells.Item(14, 3).Value = 0.1631492743891265
$ws.Cells.Item(14, 4).Value = 0.7730453408521514
$ws.Cells.Item(14, 5).Value = 0.5318168001700766
$ws.Cells.Item(14, 6).Value = 0.1702559739351273
$ws.Cells.Item(14, 7).Value = 0.7460638880729675
$ws.Cells.Item(14, 8).Value = 0.1539419144392014
$ws.Cells.Item(14, 9).Value = 0.4674183130264282

$ws.Cells.Item(15, 1).Value = "model_11_8_13"
$ws.Cells.Item(15, 2).Value = 0.8464883545132103
$ws.Cells.Item(15, 3).Value = 0.1580157920198006
$ws.Cells.Item(15, 4).Value = 0.8391915319035141
$ws.Cells.Item(15, 5).Value = 0.5505382312454239
$ws.Cells.Item(15, 6).Value = 0.1698921471834183
$ws.Cells.Item(15, 7).Value = 0.7506403923034668
$ws.Cells.Item(15, 8).Value = 0.1090753749012947
$ws.Cells.Item(15, 9).Value = 0.4487273991107941

$ws.Cells.Item(16, 1).Value = "model_11_8_19"
$ws.Cells.Item(16, 2).Value = 0.846518401184198
$ws.Cells.Item(16, 3).Value = 0.1855233636198523
$ws.Cells.Item(16, 4).Value = 0.7362762534235405
$ws.Cells.Item(16, 5).Value = 0.5306383071249554
$ws.Cells.Item(16, 6).Value = 0.1698589026927948
$ws.Cells.Item(16, 7).Value = 0.7261170744895935
$ws.Cells.Item(16, 8).Value = 0.1788821518421173
$ws.Cells.Item(16, 9).Value = 0.4685948193073273

$ws.Cells.Item(17, 1).Value = "model_11_8_16"
$ws.Cells.Item(17, 2).Value = 0.8465472481785793
$ws.Cells.Item(17, 3).Value = 0.171162958111047
$ws.Cells.Item(17, 4).Value = 0.7754958724519715
$ws.Cells.Item(17, 5).Value = 0.5363887584567777
$ws.Cells.Item(17, 6).Value = 0.1698269844055176
$ws.Cells.Item(17, 7).Value = 0.7389194965362549
$ws.Cells.Item(17, 8).Value = 0.1522797346115112
$ws.Cells.Item(17, 9).Value = 0.4628537893295288

$ws.Cells.Item(18, 1).Value = "model_11_8_2"
$ws.Cells.Item(18, 2).Value = 0.8477823776783717
$ws.Cells.Item(18, 3).Value = 0.6184859335983024
$ws.Cells.Item(18, 4).Value = 0.7314836707722451
$ws.Cells.Item(18, 5).Value = 0.7337893167637239
$ws.Cells.Item(18, 6).Value = 0.1684600561857224
$ws.Cells.Item(18, 7).Value = 0.3401249647140503
$ws.Cells.Item(18, 8).Value = 0.1821329295635223
$ws.Cells.Item(18, 9).Value = 0.265775740146637

$ws.Cells.Item(19, 1).Value = "model_11_8_18"
$ws.Cells.Item(19, 2).Value = 0.8481470644353901
$ws.Cells.Item(19, 3).Value = 0.1820664890068254
$ws.Cells.Item(19, 4).Value = 0.7643885981481919
$ws.Cells.Item(19, 5).Value = 0.537992280079984
$ws.Cells.Item(19, 6).Value = 0.1680564433336258
$ws.Cells.Item(19, 7).Value = 0.7291988730430603
$ws.Cells.Item(19, 8).Value = 0.1598137319087982
$ws.Cells.Item(19, 9).Value = 0.4612528681755066

$ws.Cells.Item(20, 1).Value = "model_11_8_5"
$ws.Cells.Item(20, 2).Value = 0.8491176555221598
$ws.Cells.Item(20, 3).Value = 0.4516460556811732
$ws.Cells.Item(20, 4).Value = 0.7657221725007248
$ws.Cells.Item(20, 5).Value = 0.6658623687218954
$ws.Cells.Item(20, 6).Value = 0.1669823080301285
$ws.Cells.Item(20, 7).Value = 0.4888650178909302
$ws.Cells.Item(20, 8).Value = 0.1589091718196869
$ws.Cells.Item(20, 9).Value = 0.3335916996002197

$ws.Cells.Item(21, 1).Value = "model_11_8_24"
$ws.Cells.Item(21, 2).Value = 0.8495027901224026
$ws.Cells.Item(21, 3).Value = 0.1997265491147611
$ws.Cells.Item(21, 4).Value = 0.7435663734660131
$ws.Cells.Item(21, 5).Value = 0.5396836955837694
$ws.Cells.Item(21, 6).Value = 0.1665560752153397
$ws.Cells.Item(21, 7).Value = 0.7134547829627991
$ws.Cells.Item(21, 8).Value = 0.1739373058080673
$ws.Cells.Item(21, 9).Value = 0.459564208984375

$ws.Cells.Item(22, 1).Value = "model_11_8_14"
$ws.Cells.Item(22, 2).Value = 0.8496696871366113
$ws.Cells.Item(22, 3).Value = 0.1831477082062988
$ws.Cells.Item(22, 4).Value = 0.8181125086850477
$ws.Cells.Item(22, 5).Value = 0.5556798613725602
$ws.Cells.Item(22, 6).Value = 0.1663713455200195
$ws.Cells.Item(22, 7).Value = 0.7282350063323975
$ws.Cells.Item(22, 8).Value = 0.1233731359243393
$ws.Cells.Item(22, 9).Value = 0.443594217300415

$ws.Cells.Item(23, 1).Value = "model_11_8_22"
$ws.Cells.Item(23, 2).Value = 0.8499280142711018
$ws.Cells.Item(23, 3).Value = 0.1725660198205951
$ws.Cells.Item(23, 4).Value = 0.790303739387698
$ws.Cells.Item(23, 5).Value = 0.5417864486566812
$ws.Cells.Item(23, 6).Value = 0.1660854667425156
$ws.Cells.Item(23, 7).Value = 0.7376687526702881
$ws.Cells.Item(23, 8).Value = 0.1422356516122818
$ws.Cells.Item(23, 9).Value = 0.4574649035930634

$ws.Cells.Item(24, 1).Value = "model_11_8_10"
$ws.Cells.Item(24, 2).Value = 0.8517972421453807
$ws.Cells.Item(24, 3).Value = 0.244927099460727
$ws.Cells.Item(24, 4).Value = 0.8124627921389396
$ws.Cells.Item(24, 5).Value = 0.5830798130105601
$ws.Cells.Item(24, 6).Value = 0.1640167683362961
$ws.Cells.Item(24, 7).Value = 0.673157811164856
$ws.Cells.Item(24, 8).Value = 0.1272053122520447
$ws.Cells.Item(24, 9).Value = 0.4162389934062958

$ws.Cells.Item(25, 1).Value = "model_11_8_11"
$ws.Cells.Item(25, 2).Value = 0.8520818658306739
$ws.Cells.Item(25, 3).Value = 0.2436541755266722
$ws.Cells.Item(25, 4).Value = 0.8161011925138925
$ws.Cells.Item(25, 5).Value = 0.5836413405070158
$ws.Cells.Item(25, 6).Value = 0.1637018024921417
$ws.Cells.Item(25, 7).Value = 0.6742926239967346
$ws.Cells.Item(25, 8).Value = 0.1247373968362808
$ws.Cells.Item(25, 9).Value = 0.4156784117221832

$ws.Cells.Item(26, 1).Value = "model_11_8_12"
$ws.Cells.Item(26, 2).Value = 0.8542397263048571
$ws.Cells.Item(26, 3).Value = 0.2305043609328755
$ws.Cells.Item(26, 4).Value = 0.8365634217388151
$ws.Cells.Item(26, 5).Value = 0.5839669253574467
$ws.Cells.Item(26, 6).Value = 0.1613136678934097
$ws.Cells.Item(26, 7).Value = 0.6860158443450928
$ws.Cells.Item(26, 8).Value = 0.1108580008149147
$ws.Cells.Item(26, 9).Value = 0.4153533577919006

